$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65

# Row 3 updates
$ws.Range("G3").Value = 1.48
$ws.Range("I3").Value = 6.25
$ws.Range("K3").Value = 2.38
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 13
$ws.Range("O3").Value = 1.22
$ws.Range("P3").Value = 4
$ws.Range("Q3").Value = 1.75
$ws.Range("R3").Value = 2.05
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 7
$ws.Range("Z3").Value = 10
$ws.Range("AA3").Value = 12
$ws.Range("AB3").Value = 26
$ws.Range("AC3").Value = 12
$ws.Range("AD3").Value = 9
$ws.Range("AJ3").Value = 19
$ws.Range("AV3").Value = 51
$ws.Range("AZ3").Value = 126

# Row 6 updates
$ws.Range("N6").Value = 13
